# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same layout as the
# other quarterly sheets) right before the existing "总计" (totals) sheet,
# and adds a corresponding summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4"
#    (i.e. right before "总计", which keeps it last).
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "2022-Q1"

# Copy header-row (B1:H1) formatting from the template quarter sheet, then
# overwrite the header captions for the new sheet.
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Copy the bold/centered/bordered "index" column style (column A) used on
# every data row, broadcasting it down to all 16 data rows at once.
$template.Range("A2").Copy($newSheet.Range("A2:A17"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Fund holdings for 2022-Q1: index, code, name, scale, stock-position-total,
# position-ratio, held-value(亿元), position-rank
$fundData = @(
    @(0, "501011", "汇添富中证中药指数（LOF）A", "24.13", "94.57", "4.52", "1.0907", 5),
    @(1, "005505", "前海开源中药研究精选股票A", "11.44", "91.95", "8.21", "0.9392", 4),
    @(2, "501012", "汇添富中证中药指数（LOF）C", "8.91", "94.57", "4.52", "0.4027", 5),
    @(3, "005506", "前海开源中药研究精选股票C", "4.62", "91.95", "8.21", "0.3793", 4),
    @(4, "673110", "西部利得新润灵活配置混合", "5.10", "81.49", "7.30", "0.3723", 1),
    @(5, "217024", "招商安盈债券", "35.05", "20.20", "0.65", "0.2278", 10),
    @(6, "005433", "申万菱信医药先锋股票", "2.20", "90.81", "5.07", "0.1115", 3),
    @(7, "519673", "银河康乐股票", "2.31", "92.35", "4.08", "0.0942", 9),
    @(8, "011383", "富安达医药创新混合", "1.68", "83.50", "3.25", "0.0546", 9),
    @(9, "005043", "国寿安保健康科学混合A", "0.99", "85.72", "2.85", "0.0282", 7),
    @(10, "005044", "国寿安保健康科学混合C", "0.87", "85.72", "2.85", "0.0248", 7),
    @(11, "006478", "长盛多因子策略优选股票", "0.51", "84.41", "4.83", "0.0246", 1),
    @(12, "004351", "汇丰晋信珠三角区域发展混合", "0.51", "93.92", "4.63", "0.0236", 8),
    @(13, "001861", "富安达健康人生灵活配置混合", "0.61", "82.18", "2.97", "0.0181", 10),
    @(14, "009502", "国寿安保创新医药股票A", "0.54", "81.60", "2.92", "0.0158", 6),
    @(15, "009503", "国寿安保创新医药股票C", "0.20", "81.60", "2.92", "0.0058", 6)
)

$r = 2
foreach ($row in $fundData) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q1 and keep the index column (A) sequential.
# ---------------------------------------------------------------------

# Pre-format the index column for the (now) 6 data rows using the existing
# styled cell as a template, so the new row matches the others.
$totalSheet.Range("A2").Copy($totalSheet.Range("A2:A7"))

$totalsData = @(
    @(0, "2022-Q1", 16, 3.81),
    @(1, "2021-Q4", 7, 1.94),
    @(2, "2021-Q3", 11, 1.62),
    @(3, "2021-Q2", 11, 1.04),
    @(4, "2021-Q1", 5, 0.35),
    @(5, "2020-Q4", 2, 0.13)
)

$r = 2
foreach ($row in $totalsData) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
